$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 data with the new person's details
$ws.Range("A2").Value = "Rico Putra Pradana"
$ws.Range("B2").Value = "rico.putra@outlook.co.id"
$ws.Range("C2").Value = "22 Tahun 7 Bulan"

# Resize columns A and B
# (Excel's ColumnWidth setter rounds to the nearest pixel, i.e. 1/6 of a
# character width, before it is persisted as the sheet's <col width=.../>.
# Subtracting the standard 5-pixel cell-padding offset of 5/6 lands us in
# the pixel bucket closest to the target stored width.)
$ws.Columns.Item(1).ColumnWidth = 18.877604166666668
$ws.Columns.Item(2).ColumnWidth = 24.877604166666668
